$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.726.64"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "3.518.16"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'584.98"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'175.98"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.513.84"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("E10").Value = "  -2.73%  "
$ws.Range("D11").Value = "'6.91"
$ws.Range("E12").Value = "  -3.13%  "
$ws.Range("D13").Value = "4.128.94"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "'30.60"
$ws.Range("E14").Value = "  -5.16%  "
$ws.Range("D15").Value = "'0.133"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "66.717.04"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").Value = "3.523.58"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").Value = "'14.00"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "'382.70"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").Value = "'7.93"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "'0.551"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'72.35"
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("E28").Value = "  -4.83%  "
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'24.64"
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("D32").Value = "'5.93"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("E34").Value = "  -5.68%  "
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D38").Value = "'30.11"
$ws.Range("E38").Value = "  +13.15%  "
$ws.Range("D39").Value = "'160.73"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").Value = "'0.898"
$ws.Range("E41").Value = "  -5.26%  "
$ws.Range("D43").Value = "'6.49"
$ws.Range("E43").Value = "  -5.07%  "
$ws.Range("D44").Value = "'2.55"
$ws.Range("E44").Value = "  -10.02%  "
$ws.Range("D45").Value = "2.726.62"
$ws.Range("E45").Value = "  -4.72%  "
$ws.Range("D46").Value = "'0.0703"
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value = "'40.79"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").Value = "'25.05"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").Value = "'324.23"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("E51").Value = "  -3.85%  "
